$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Re-label task IDs: single-digit "tN" (t1..t9) become zero-padded "t0N"
#    (t10 / t11 are already two digits and stay as-is).
for ($r = 2; $r -le 12; $r++) {
    $idCell = $ws.Cells.Item($r, 1)
    $id = $idCell.Value2
    if ($id -match '^t(\d)$') {
        $idCell.Value2 = "t0" + $matches[1]
    }
}

# 2) Every task is now marked done.
$ws.Range("C2:C12").Value2 = "yes"

# 3) The Queue column values are no longer used - clear them out.
$ws.Range("F2:F5").ClearContents()

# 4) Re-sort the table by ID (column A) ascending instead of by Queue.
$lo = $ws.ListObjects.Item(1)
$lo.Sort.SortFields.Clear()
$lo.Sort.SortFields.Add($ws.Range("A1:A12"))
$lo.Sort.Apply()

# 5) The long wrapped-text row (previously row 2) now lands on row 12 -
#    let row heights follow the content back to the default, then restore
#    the taller height on the row that now holds the long text.
$ws.Range("A1:A12").EntireRow.AutoFit()
$ws.Rows.Item(12).RowHeight = 28.8

# 6) Update the saved selection.
$ws.Range("B18").Select()
